$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XLFormula")
$ws.Activate()

# Append the new "neg test case" sample rows (B46:B57) below the existing data.
$values = @(1, 1, 3, 3, 50, 4, 3, 1, 3, 1, 1, 1)
$startRow = 46
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Scroll the view down and leave the selection on the last entered cell,
# matching the author's final on-screen state.
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B54").Select()
